{"js": "// no-op test\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# no-op test\n$d = $word.ActiveDocument\nWrite-Output $d.Paragraphs.Count\n"}
